$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up the existing EB94.3 row (currently row 3, "ZA7780") ---
# The timeframe text had a stray " - " instead of "-", and the description
# (COVID-19 Pandemic) was missing entirely.
$ws.Cells.Item(3, 3).Value = "February-March 2021"
$ws.Cells.Item(3, 4).Value = "COVID-19 Pandemic"

# --- Insert the new EB95.3 row right under the header row ---
$ws.Rows.Item(2).Insert()
$ws.Cells.Item(2, 1).Value = "ZA7783"
$ws.Cells.Item(2, 2).Value = "'95.3"
$ws.Cells.Item(2, 3).Value = "June-July 2021"
$ws.Cells.Item(2, 4).Value = "COVID-19 Pandemic"

# --- Column widths for B (wave) and C (timeframe) ---
$ws.Columns.Item(2).ColumnWidth = 12.666666666666668
$ws.Columns.Item(3).ColumnWidth = 25.5

# --- Selection moves to D3 ---
$ws.Range("D3").Select()
